$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = "asdasd"
$ws.Range("F8").Value = "asdasd"
$ws.Range("I10").Value = "asdasd"

$ws.Range("U4").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
